$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 10 data values (columns A & B, rows 1-5) to the new case-1 data.
$ws.Range("A1").Value = 0.025319509551944442
$ws.Range("B1").Value = -0.025319510267234623
$ws.Range("A2").Value = -0.040594747259237472
$ws.Range("B2").Value = 0.040594746539864057
$ws.Range("A3").Value = -0.014060498836495407
$ws.Range("B3").Value = 0.014060498077372313
$ws.Range("A4").Value = 0.057646366125184395
$ws.Range("B4").Value = -0.057646366855596404
$ws.Range("A5").Value = -0.015368927939619836
$ws.Range("B5").Value = 0.015368927170522332

# Columns A & B both end up at the same (narrower) width after the data
# refresh - closest value reachable through the ColumnWidth property's
# pixel-grid rounding to the target 14.42578125 characters.
$ws.Columns.Item(1).ColumnWidth = 13.666666666666666
$ws.Columns.Item(2).ColumnWidth = 13.666666666666666
